$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '65.516.45'
$ws.Range("E2").Value = '  -3.33%  '

$ws.Range("D3").Value = '3.499.10'
$ws.Range("E3").Value = '  -0.09%  '

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '1.00'
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = '  +0.08%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '555.70'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  -0.18%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '179.52'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  -6.13%  '

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.638'
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  +3.88%  '

$ws.Range("E8").Value = '  -0.09%  '

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.631'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  -1.42%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.154'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  +2.66%  '

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '53.83'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  -6.57%  '

$ws.Range("E12").Value = '  -1.34%  '

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '9.26'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  -2.44%  '

$ws.Range("D14").Value = '4.059.18'
$ws.Range("E14").Value = '  -0.27%  '

$ws.Range("D15").Value = '3.513.47'
$ws.Range("E15").Value = '  +0.24%  '

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '18.45'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  +0.43%  '

$ws.Range("E17").Value = '  +0.25%  '

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '12.14'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  +2.21%  '

$ws.Range("D19").Value = '65.631.02'
$ws.Range("E19").Value = '  -3.63%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '0.996'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  -1.34%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '416.54'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  +2.69%  '

$ws.Range("E22").Value = '  +2.61%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '86.12'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  +1.41%  '

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '4.12'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  -2.20%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '12.87'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  +8.32%  '

$ws.Range("E26").Value = '  -10.10%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '2.86'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  -2.27%  '

$ws.Range("E28").Value = '  -3.96%  '

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '9.06'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  +4.69%  '

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '30.36'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  -0.30%  '

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '6.49'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  -5.96%  '

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '609.49'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  -12.38%  '

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '11.69'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  -0.41%  '

$ws.Range("E34").Value = '  -1.32%  '

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '59.76'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  -1.99%  '

$ws.Range("E36").Value = '  +8.90%  '

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '37.35'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  -4.31%  '

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.995'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  -0.38%  '

$ws.Range("D39").Value = '3.410.05'
$ws.Range("E39").Value = '  +11.29%  '

$ws.Range("D40").Value = '0.0₃0793'
$ws.Range("E40").Value = '  -4.53%  '

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.380'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  -5.97%  '

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '1.00'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  +0.23%  '

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '3.25'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  -4.33%  '

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '2.85'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  -4.15%  '

$ws.Range("E45").Value = '  -9.55%  '

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.0415'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  -1.68%  '

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '3.25'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  +1.27%  '

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '2.71'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  -2.88%  '

$ws.Range("E49").Value = '  +1.60%  '

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '8.45'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  -4.06%  '

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '137.62'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  -2.21%  '
